$d = $word.ActiveDocument

# wdReplaceAll = 2 ; wdFindContinue = 1
$wdReplaceAll = 2
$wdFindContinue = 1

# 1) "Case No:" FILLIN / _Q1_CASENUM  -> merge split runs (removes the
#    stray <w:proofErr w:type="gramEnd"/> Word had inserted between them)
$d.Content.Find.Execute("FILLIN_Q1_CASENUM", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "FILLIN_Q1_CASENUM", $wdReplaceAll) | Out-Null

# 2) "FILLIN_" + "PARAGRAPH_GENERAL" -> single run "FILLIN_PARAGRAPH_GENERAL"
$d.Content.Find.Execute("FILLIN_PARAGRAPH_GENERAL", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "FILLIN_PARAGRAPH_GENERAL", $wdReplaceAll) | Out-Null

# 3) "FILLIN_" + " PARAGRAPH_CONFIDENTIALITY" -> "FILLIN_PARAGRAPH_CONFIDENTIALITY"
#    (also drops the extra inner space that used to separate the two runs)
$d.Content.Find.Execute("FILLIN_ PARAGRAPH_CONFIDENTIALITY", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "FILLIN_PARAGRAPH_CONFIDENTIALITY", $wdReplaceAll) | Out-Null

# 4) " FILLIN_Q12" -> " FILLIN_PARAGRAPH_REVIEW_OF_RECORDS"
#    MatchWholeWord so "FILLIN_Q125" (a different placeholder, further down
#    the document) is left untouched.
$d.Content.Find.Execute(" FILLIN_Q12", $true, $true, $false, $false, $false, $true, $wdFindContinue, $false, " FILLIN_PARAGRAPH_REVIEW_OF_RECORDS", $wdReplaceAll) | Out-Null

# 5) "FILLIN_PARAGRAPH_ATTITUDE_&_BEHAVIOR" -> "FILLIN_PARAGRAPH_ATTITUDE_AND_BEHAVIOR"
$d.Content.Find.Execute("FILLIN_PARAGRAPH_ATTITUDE_&_BEHAVIOR", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "FILLIN_PARAGRAPH_ATTITUDE_AND_BEHAVIOR", $wdReplaceAll) | Out-Null

# 6) "FILLIN_" + "PARAGRAPH_INTELLECTUAL_" + "FUNCTIONING" (x2) -> single run
$d.Content.Find.Execute("FILLIN_PARAGRAPH_INTELLECTUAL_FUNCTIONING", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "FILLIN_PARAGRAPH_INTELLECTUAL_FUNCTIONING", $wdReplaceAll) | Out-Null

# 7) "FILLIN_" + "PARAGRAPH_FUND_OF_KNOWLEDGE_INFORMATION" (x2) -> single run
$d.Content.Find.Execute("FILLIN_PARAGRAPH_FUND_OF_KNOWLEDGE_INFORMATION", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "FILLIN_PARAGRAPH_FUND_OF_KNOWLEDGE_INFORMATION", $wdReplaceAll) | Out-Null

# 8) "the  FILLIN" (double space) + "_Q125 range." -> "the FILLIN_Q125 range." (single space, merged run)
$d.Content.Find.Execute("the  FILLIN_Q125 range.", $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, "the FILLIN_Q125 range.", $wdReplaceAll) | Out-Null
